$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerMap = @{
    "D1"  = "이슬점(Td)"
    "E1"  = "증산(HD)"
    "F1"  = "대기압"
    "G1"  = "PWS(?)"
    "H1"  = "수분량"
    "I1"  = "엔탈피"
    "J1"  = "비부피"
    "K1"  = "밀도"
    "L1"  = "PW(?)수증기"
    "M1"  = "절대습도"
    "N1"  = "절대AH(?)"
    "O1"  = "실내온도표준편차"
    "P1"  = "실내습도표준편차"
    "Q1"  = "이슬점(Td)표준편차"
    "R1"  = "증산(HD)표준편차"
    "S1"  = "대기압표준편차"
    "T1"  = "PWS(?)표준편차"
    "U1"  = "수분량표준편차"
    "V1"  = "엔탈피표준편차"
    "W1"  = "비부피표준편차"
    "X1"  = "밀도표준편차"
    "Y1"  = "PW(?)수증기표준편차"
    "Z1"  = "절대습도표준편차"
    "AA1" = "절대AH(?)표준편차"
    "AB1" = "최소실내온도"
    "AC1" = "최소실내습도"
    "AD1" = "최소이슬점(Td)"
    "AE1" = "최소증산(HD)"
    "AF1" = "최소대기압"
    "AG1" = "최소PWS(?)"
    "AH1" = "최소수분량"
    "AI1" = "최소엔탈피"
    "AJ1" = "최소비부피"
    "AK1" = "최소밀도"
    "AL1" = "최소PW(?)수증기"
    "AM1" = "최소절대습도"
    "AN1" = "최소절대AH(?)"
    "AO1" = "최대실내온도"
    "AP1" = "최대실내습도"
    "AQ1" = "최대이슬점(Td)"
    "AR1" = "최대증산(HD)"
    "AS1" = "최대대기압"
    "AT1" = "최대PWS(?)"
    "AU1" = "최대수분량"
    "AV1" = "최대엔탈피"
    "AW1" = "최대비부피"
    "AX1" = "최대밀도"
    "AY1" = "최대PW(?)수증기"
    "AZ1" = "최대절대습도"
    "BA1" = "최대절대AH(?)"
    "BB1" = "오전적정증산(HD)누적시간"
    "BC1" = "일출일몰적정증산(HD)누적시간"
    "BD1" = "하루적정증산(HD)누적시간"
    "BE1" = "오후심각증산(HD)누적시간"
    "BF1" = "오전심각증산(HD)누적시간"
    "BG1" = "12도이하온도누적시간"
    "BH1" = "30도이상온도누적시간"
    "BI1" = "적정온도누적시간"
    "BJ1" = "주간평균온도"
    "BK1" = "야간평균온도"
    "BL1" = "오후부터일몰까지평균온도"
    "BM1" = "적정습도누적시간"
    "BN1" = "일출전후한시간평균온도"
    "BO1" = "일몰전후한시간평균온도"
    "BP1" = "적정온도변화폭(하위)"
    "BQ1" = "적정온도변화폭(상위)"
}

foreach ($cellRef in $headerMap.Keys) {
    $ws.Range($cellRef).Value = $headerMap[$cellRef]
}
